# Se importo nuevos archivos
# Fix column alignment from a bad import: birth-date values that landed in
# column C actually belong in column D, and "Carrera" values that landed in
# column H actually belong in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # C
    $dCell = $ws.Cells.Item($r, 4)   # D
    $gCell = $ws.Cells.Item($r, 7)   # G
    $hCell = $ws.Cells.Item($r, 8)   # H

    $cText = $cCell.Text
    if ($cText -ne $null -and $cText -ne "") {
        $cCell.ClearContents()
        $dCell.Value2 = "'" + $cText
        $dCell.Style = "Normal"
    }

    $hText = $hCell.Text
    if ($hText -ne $null -and $hText -ne "") {
        $hCell.ClearContents()
        $gCell.Value2 = "'" + $hText
        $gCell.Style = "Normal"
    }
}

# Update the active selection to match the saved view state.
$ws.Range("A2:H12").Select()
$excel.ActiveCell = $ws.Range("A2")
